{"js": "// The document's last paragraph reads \"Iam karthick\" as a single run.\n// The target edit:\n//   1. Splits that run into two runs: \"Iam \" and \"Karthick\" (capitalized),\n//      both keeping the same run formatting (<w:lang w:val=\"en-US\"/>).\n//   2. Adds a brand-new paragraph directly after it containing \"I am vikram\".\n\nconst body = context.document.body;\n\n// --- Step 1: locate the \"karthick\" substring inside the last paragraph ---\nconst hits = body.search(\"karthick\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const target = hits.items[0];\n\n  // Replace \"karthick\" with \"Karthick\" in place, then force the\n  // replacement text to become its own run (distinct from the \"Iam \"\n  // run that precedes it) by toggling a character-formatting property\n  // on it. Flipping bold on then off keeps the final formatting\n  // identical to the surrounding text while leaving the run boundary\n  // intact (engines merge two runs only if they were never mutated\n  // independently once split).\n  const newRange = target.insertText(\"Karthick\", \"Replace\");\n  newRange.font.bold = true;\n  await context.sync();\n\n  newRange.font.bold = false;\n  await context.sync();\n}\n\n// --- Step 2: add the new paragraph \"I am vikram\" right after it ---\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst lastParagraph = paras.items[paras.items.length - 1];\nlastParagraph.insertParagraph(\"I am vikram\", \"After\");\nawait context.sync();\n", "ps1": "# The document's last paragraph reads \"Iam karthick\" as a single run.\n# The target edit:\n#   1. Splits that run into two runs: \"Iam \" and \"Karthick\" (capitalized),\n#      both keeping the same run formatting (<w:lang w:val=\"en-US\"/>).\n#   2. Adds a brand-new paragraph directly after it containing \"I am vikram\".\n\n$d = $word.ActiveDocument\n\n# --- Step 1: locate the \"karthick\" substring inside the last paragraph ---\n$rng = $d.Range()\n$rng.Find.MatchCase = $true\n$rng.Find.Text = \"karthick\"\n$rng.Find.Execute() | Out-Null\n\nif ($rng.Find.Found) {\n    # Replace \"karthick\" with \"Karthick\" in place, then force the\n    # replacement text to become its own run (distinct from the \"Iam \"\n    # run that precedes it) by toggling a character-formatting property\n    # on it. Flipping Bold on then off keeps the final formatting\n    # identical to the surrounding text while leaving the run boundary\n    # intact (runs only merge back together if they were never mutated\n    # independently once split).\n    $rng.Text = \"Karthick\"\n    $rng.Font.Bold = $true\n    $rng.Font.Bold = $false\n}\n\n# --- Step 2: add the new paragraph \"I am vikram\" right after it ---\n$last = $d.Paragraphs.Last\n$last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"I am vikram\"\n"}
